$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 35
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 35
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 80
$ws.Range("N5").ClearContents()

$ws.Range("H9").Value = 225.63637
$ws.Range("I9").Value = 97.42856999999999
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 97.42856999999999
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = 71.57143000000001
$ws.Range("N9").Value = -788

$ws.Range("H28").Value = 1063.6
$ws.Range("I28").Value = 897.4
$ws.Range("J28").Value = 1396
$ws.Range("K28").Value = 897.4
$ws.Range("L28").Value = 1396
$ws.Range("M28").Value = -412.4
$ws.Range("N28").Value = -2366

$ws.Range("H40").Value = 1044.7894
$ws.Range("J40").Value = 1020.0571
$ws.Range("L40").Value = 1020.0571
$ws.Range("N40").Value = -1370.0571

$ws.Range("H74").Value = 3919.5293
$ws.Range("I74").Value = 3825
$ws.Range("J74").Value = 4003.5557
$ws.Range("K74").Value = 3825
$ws.Range("L74").Value = 4003.5557
$ws.Range("M74").Value = -2889
$ws.Range("N74").Value = -5875.5557

$ws.Range("H76").Value = 3672
$ws.Range("I76").Value = 3133.3333
$ws.Range("J76").Value = 3902.8572
$ws.Range("K76").Value = 3133.3333
$ws.Range("L76").Value = 3902.8572
$ws.Range("M76").Value = -2818.3333
$ws.Range("N76").Value = -4532.8572

$ws.Range("H77").Value = 3919.5293
$ws.Range("I77").Value = 3825
$ws.Range("J77").Value = 4003.5557
$ws.Range("K77").Value = 19125
$ws.Range("L77").Value = 20017.7785
$ws.Range("M77").Value = -14445
$ws.Range("N77").Value = -29377.7785

$ws.Range("H79").Value = 3672
$ws.Range("I79").Value = 3133.3333
$ws.Range("J79").Value = 3902.8572
$ws.Range("K79").Value = 3133.3333
$ws.Range("L79").Value = 3902.8572
$ws.Range("M79").Value = -2041.3333
$ws.Range("N79").Value = -6086.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2237.2593
$ws.Range("I2").Value = 2154.2273
$ws.Range("K2").Value = 2154.2273
$ws.Range("M2").Value = -2041.2273

$ws.Range("H116").Value = 2237.2593
$ws.Range("I116").Value = 2154.2273
$ws.Range("K116").Value = 2154.2273
$ws.Range("M116").Value = 139.7727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2237.2593
$ws.Range("I3").Value = 2154.2273
$ws.Range("K3").Value = 2154.2273
$ws.Range("M3").Value = -2040.2273

$ws.Range("H105").Value = 166668800
$ws.Range("I105").Value = 166668800
$ws.Range("K105").Value = 166668800
$ws.Range("M105").Value = -166667053

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 18520332
$ws.Range("I16").Value = 1937.75
$ws.Range("J16").Value = 33335048
$ws.Range("K16").Value = 1937.75
$ws.Range("L16").Value = 33335048
$ws.Range("M16").Value = -1650.75
$ws.Range("N16").Value = -33335622

$ws.Range("H50").Value = 21909.715
$ws.Range("J50").Value = 21909.715
$ws.Range("L50").Value = 21909.715
$ws.Range("N50").Value = -23159.715

$ws.Range("H113").Value = 18520332
$ws.Range("I113").Value = 1937.75
$ws.Range("J113").Value = 33335048
$ws.Range("K113").Value = 1937.75
$ws.Range("L113").Value = 33335048
$ws.Range("M113").Value = 232.25
$ws.Range("N113").Value = -33339388

$ws.Range("H122").Value = 1731.3334
$ws.Range("I122").Value = 1352.6
$ws.Range("J122").Value = 3625
$ws.Range("K122").Value = 4057.8
$ws.Range("L122").Value = 10875
$ws.Range("M122").Value = -1607.8
$ws.Range("N122").Value = -15775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 174.93333
$ws.Range("I2").Value = 242
$ws.Range("J2").Value = 74.333336
$ws.Range("K2").Value = 1452
$ws.Range("L2").Value = 446.000016
$ws.Range("M2").Value = -1339
$ws.Range("N2").Value = -672.000016

$ws.Range("H92").Value = 916
$ws.Range("I92").Value = 1066.6666
$ws.Range("J92").Value = 690
$ws.Range("K92").Value = 3199.9998
$ws.Range("L92").Value = 2070
$ws.Range("M92").Value = -1951.9998
$ws.Range("N92").Value = -4566

$ws.Range("H131").Value = 880.4918
$ws.Range("J131").Value = 887.7966300000001
$ws.Range("L131").Value = 2663.38989
$ws.Range("N131").Value = -12743.38989

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 36912.47
$ws.Range("I70").Value = 51027.227
$ws.Range("J70").Value = 5860
$ws.Range("K70").Value = 51027.227
$ws.Range("L70").Value = 5860
$ws.Range("M70").Value = -50757.227
$ws.Range("N70").Value = -6400

$ws.Range("H73").Value = 36912.47
$ws.Range("I73").Value = 51027.227
$ws.Range("J73").Value = 5860
$ws.Range("K73").Value = 51027.227
$ws.Range("L73").Value = 5860
$ws.Range("M73").Value = -50091.227
$ws.Range("N73").Value = -7732

$ws.Range("H132").Value = 75882.3
$ws.Range("I132").Value = 68234
$ws.Range("J132").Value = 85442.664
$ws.Range("K132").Value = 204702
$ws.Range("L132").Value = 256327.992
$ws.Range("M132").Value = -202172
$ws.Range("N132").Value = -261387.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H61").Value = 2612
$ws.Range("I61").Value = 1779.2
$ws.Range("K61").Value = 1779.2
$ws.Range("M61").Value = -1577.2

$ws.Range("H113").Value = 2612
$ws.Range("I113").Value = 1779.2
$ws.Range("K113").Value = 1779.2
$ws.Range("M113").Value = 390.8

$ws.Range("H122").Value = 3446.6985
$ws.Range("I122").Value = 3007
$ws.Range("J122").Value = 3596.383
$ws.Range("K122").Value = 9021
$ws.Range("L122").Value = 10789.149
$ws.Range("M122").Value = -6571
$ws.Range("N122").Value = -15689.149

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 4015
$ws.Range("I21").Value = 4015
$ws.Range("K21").Value = 4015
$ws.Range("M21").Value = -3780

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()

$ws.Range("H29").Value = 4000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -4580

$ws.Range("H35").Value = 4015
$ws.Range("I35").Value = 4015
$ws.Range("K35").Value = 4015
$ws.Range("M35").Value = -3725

$ws.Range("H113").Value = 1098.5714
$ws.Range("I113").Value = 1098.5714
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3295.7142
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1125.7142
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 1970.1
$ws.Range("I126").Value = 1966.7778
$ws.Range("K126").Value = 5900.3334
$ws.Range("M126").Value = -3430.3334
